# Apply the authored change:
#  - Sheet "PIE" (sheet1): no data change, only the saved selection moves to B7:E7.
#  - Sheet "BAR" (sheet2): the blank row 1 is deleted, shifting every data row
#    up by one (old row 2 header -> new row 1, ..., old row 29 -> new row 28),
#    and the saved selection moves to L29. BAR stays the active/selected tab,
#    matching the source workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "PIE" : selection only -------------------------------------------------
$wsPie = $wb.Worksheets.Item("PIE")
$wsPie.Range("B7:E7").Select()

# --- Sheet "BAR" : delete the leading blank row, then fix the selection ----------
$wsBar = $wb.Worksheets.Item("BAR")
$wsBar.Rows.Item(1).Delete()

# Keep/restore BAR as the active sheet & tab (selecting on PIE above would
# otherwise have made PIE the active tab), then set its saved selection.
$wsBar.Activate()
$wsBar.Range("L29").Select()
